$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualização de versões instaladas: Word, Power Point e Project passaram
# de 2019 para 2021 (ficando "Atualizado" na coluna de status, que é
# calculada por fórmula e recalcula automaticamente).
$ws.Range("B7").Value = 2021
$ws.Range("B8").Value = 2021
$ws.Range("B10").Value = 2021

# A seleção ativa da planilha passou para B18.
$ws.Range("B18").Select()
